# Insert a new "Jammu and Kashmir" state block (5 rows of data) right
# before the existing "Jharkhand" block, pushing every subsequent state's
# rows down by 5. The previously-blank trailing rows (which only held a
# carried-over date style) are removed so the sheet stays the same total
# length (header + 29 states * 5 rows = 146 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Jharkhand" (and everything after it) currently starts at row 47.
# Insert 5 fresh rows there for the new state.
$ws.Rows("47:51").Insert()

$newStateName = "Jammu and Kashmir"
$dates = @(44197, 44228, 44256, 44287, 44317)
$values = @(0.48, 0.49, 0.49, 0.5, 0.51)

for ($i = 0; $i -lt 5; $i++) {
    $r = 47 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $newStateName
    $ws.Cells.Item($r, 3).Value = $values[$i]
}

# The insert shifted the 5 trailing (formerly empty) rows down to
# 147:151 -- drop them so the sheet keeps its original extent (A1:C146).
$ws.Rows("147:151").Delete()

# Update the view to match where the edit was made.
$null = $ws.Range("C52").Select()
$excel.ActiveWindow.ScrollRow = 37
